$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row right under the header so the newest circular
#        (14-10-2025) becomes row 2, pushing every existing data row down
#        by one (their own Sl.no. values travel with them unchanged). ---
$ws.Rows.Item(2).Insert()

# Row-insert in this engine does not clone the neighbouring row's number
# format / alignment, so copy it explicitly from the row that used to be
# the old row 2 (now row 3) down onto the freshly inserted row 2.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Populate the new top row with the latest Nalco ingot circular. ---
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 285.05
$ws.Range("E2").Value = "14-10-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf"

# --- 3. Rebuild the hyperlinks. Inserting a row shifts the underlying
#        cell data but this engine leaves existing Hyperlink anchors
#        pinned to their old row, so drop them all and re-add them,
#        now that every row/value is in its final place. ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Adding a Hyperlink auto-applies Excel's blue/underlined "Hyperlink" cell
# style. The source sheet instead renders hyperlinked cells with the same
# plain centered style as the rest of the row, so restore that look by
# re-pasting the (already-correct) formatting from column E of each row
# onto column F.
$ws.Range("E2:E13").Copy()
$ws.Range("F2:F13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Done updating Nalco prices sheet"
